$wb = $excel.ActiveWorkbook

# "questions" sheet is the 2nd sheet in the workbook
$ws = $wb.Worksheets.Item("questions")
$ws.Activate()

# Insert 3 new rows before the current row 4 ("In this demo, I want you to say OK.")
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$ws.Range("B4").Value = "Number"
$ws.Range("B5").Value = "Date"
$ws.Range("B6").Value = "Time"

$ws.Range("A4").Value = "How many application with Amazon Connect?"
$ws.Range("A5").Value = "The date for us to call you back."
$ws.Range("A6").Value = "What time will you prefer for the callback?"

# Excel auto-fit the first column to the new (now shorter) longest entry.
# (The resulting best-fit width, in "characters", comes out to ~87.33.)
$ws.Columns.Item(1).ColumnWidth = 86.5

$ws.Range("A7").Select()

$ws1 = $wb.Worksheets.Item("configures")
$ws1.Activate()
$ws1.Range("A1:B4").Select()

$ws3 = $wb.Worksheets.Item("receivers")
$ws3.Activate()
$ws3.Range("B2").Select()
